$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension / layout happens automatically as cells are populated.

# Row 1: Bank / Country query
$ws.Range("A1").Value = "Bank"
$ws.Range("B1").Value = "Select * from dbo.country where countrypseudo ='{{ref_text}}' "

# Row 2: Stock / Stock query (keep existing content & style, just reposition)
$ws.Range("A2").Value = "Stock"
$ws.Range("B2").Value = "Select  Stock_Item_Db_Id, Stock_Item_Num, Responsible_Representative_Id, re.display_name, Representative_status_code,Stock_Item_Status_Code, Stock_Type_Code, st.Last_Updated_Datetime,st.Last_Updated_Id
from NFS_DBO.Stock_Item st (nolock)
JOIN NFS_DBO.Representative_Table re
 ON st.Responsible_representative_id = re.Representative_id
Where
Stock_Item_Num = '{{ref_text}}'
order by st.Last_Updated_Datetime"

# Row 3: Product / Product query (new row)
$ws.Range("A3").Value = "Product"
$ws.Range("B3").Value = "Select * from dbo.Products where ProductName ='{{ref_text}}' "

# Selection / view state
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("B4").Select()
